$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")
$ws.Range("A118").Formula = "=B118+C118*256"
$ws.Range("B118").Value = 2
$ws.Range("C118").Value = 7
$ws.Range("D118").Value = "AgrabahShop"
$ws.Range("E118").Value = "The Peddler" + [char]8217 + "s Shop"
$ws.Range("F118").Value = 0
$v2 = $ws.Range("E118").Text
$v2
